$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J3").Value = '病灶,中枢神经系统转移,局部治疗'
$ws.Range("J4").Value = '切片,活检'
$ws.Range("J6").Value = '胸腔积液,心包积液,洗脱期,前三个月,腹水'
$ws.Range("J8").Value = '一线治疗,KRAS,突变'
$ws.Range("J9").Value = '辅助治疗,洗脱期,6个月'
$ws.Range("J11").Value = '炎症,感染,类固醇,ILD'
$ws.Range("J12").Value = 'T-DXd,第4,方案,局部治疗,脑转移,中枢神经系统转移'
$ws.Range("J13").Value = 'MP相关,铂,毒性,免疫,禁忌症,帕博利珠,自身免疫,医学禁忌'
$ws.Range("J14").Value = '人群,HER2突变,一线治疗,突变'
$ws.Range("J15").Value = '方案,样本,肿瘤样本,突变'
$ws.Range("J16").Value = 'ctDNA,EGFR,突变,数据且,方案,ROS1、PDL1,ALK'
$ws.Range("J17").Value = '靶病灶,病灶,局部治疗,中枢神经系统转移,RECIST'
$ws.Range("J18").Value = '胸腔积液,心包积液,洗脱期,T-DXd,肺炎,肺癌,腹水'
$ws.Range("J19").Value = 'HER2,基因组,ctDNA,突变,基因突变,方案,基因,HER2突变'
$ws.Range("J21").Value = '方案,洗脱期,治疗洗脱期'
$ws.Range("J22").Value = '方案,洗脱期,治疗洗脱期'
$ws.Range("J24").Value = '洗脱期,2周'
$ws.Range("J25").Value = '胸腔积液,洗脱期,穿刺'
$ws.Range("J26").Value = '胸腔积液,洗脱期,穿刺'
$ws.Range("J27").Value = '胸腔积液,洗脱期,穿刺'
$ws.Range("J28").Value = '日本,CNS转移,方案,局部治疗,CNS,中枢神经系统转移'
$ws.Range("J34").Value = '观察性,ICF,临床试验,临床研究,第21'
$ws.Range("J36").Value = '腺癌,组织学类型,小细胞肺癌'
$ws.Range("J39").Value = '2个月,肿瘤组织,切片,样本,6个月'
$ws.Range("J47").Value = '胸腔积液,ICF,洗脱期,方案,两周'
$ws.Range("J49").Value = 'EGFR,突变,美国,莫博替尼,日本,ALK'
$ws.Range("J50").Value = 'ICF,T-DXd'
$ws.Range("J51").Value = '根治性治疗,放疗,IV期'
$ws.Range("J52").Value = '根治性治疗,放疗,IV期'
$ws.Range("J53").Value = '腺癌,组织学类型,小细胞肺癌'
$ws.Range("J54").Value = '代谢失,MSUD,感染,T-DXd'
$ws.Range("J59").Value = 'ICF,洗脱期,脑转移,中枢神经系统'
$ws.Range("J60").Value = '一线治疗,洗脱期,姑息性全身治疗,铂,全身治疗,3周,方案,局部晚期,铂类,免疫,6个月,6个月后'
$ws.Range("J61").Value = '一周,方案,前一周'
